$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5999
$ws.Range("I18").Value = 3998.5
$ws.Range("K18").Value = 3998.5
$ws.Range("M18").Value = -3714.5

$ws.Range("H74").Value = 4999
$ws.Range("I74").Value = 4999
$ws.Range("K74").Value = 4999
$ws.Range("M74").Value = -4063

$ws.Range("H77").Value = 4999
$ws.Range("I77").Value = 4999
$ws.Range("K77").Value = 24995
$ws.Range("M77").Value = -20315

$ws.Range("H112").Value = 905.5714
$ws.Range("J112").Value = 1113
$ws.Range("L112").Value = 3339
$ws.Range("N112").Value = -5555

$ws.Range("H125").Value = 7748
$ws.Range("I125").Value = 7998
$ws.Range("J125").Value = 7498
$ws.Range("K125").Value = 71982
$ws.Range("L125").Value = 67482
$ws.Range("M125").Value = -69522
$ws.Range("N125").Value = -72402

$ws.Range("H129").Value = 2398.5
$ws.Range("I129").Value = 2197
$ws.Range("J129").Value = 2427.2856
$ws.Range("K129").Value = 6591
$ws.Range("L129").Value = 7281.8568
$ws.Range("M129").Value = -1591
$ws.Range("N129").Value = -17281.8568

$ws.Range("H137").Value = 2855
$ws.Range("I137").Value = 1782.5
$ws.Range("J137").Value = 5000
$ws.Range("K137").Value = 5347.5
$ws.Range("L137").Value = 15000
$ws.Range("M137").Value = -2797.5
$ws.Range("N137").Value = -20100

$ws.Range("H138").Value = 6816.147
$ws.Range("I138").Value = 3076.3333
$ws.Range("J138").Value = 7617.5356
$ws.Range("K138").Value = 9228.999899999999
$ws.Range("L138").Value = 22852.6068
$ws.Range("M138").Value = -4088.999899999999
$ws.Range("N138").Value = -33132.6068

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2951.5833
$ws.Range("I45").Value = 2637.3333
$ws.Range("J45").Value = 3265.8333
$ws.Range("K45").Value = 2637.3333
$ws.Range("L45").Value = 3265.8333
$ws.Range("M45").Value = -2260.3333
$ws.Range("N45").Value = -4019.8333

$ws.Range("H110").Value = 2756.5625
$ws.Range("I110").Value = 3013.2144
$ws.Range("J110").Value = 960
$ws.Range("K110").Value = 3013.2144
$ws.Range("L110").Value = 960
$ws.Range("M110").Value = -968.2143999999998
$ws.Range("N110").Value = -5050

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1015.2857
$ws.Range("I94").Value = 876.75
$ws.Range("K94").Value = 876.75
$ws.Range("M94").Value = -425.75

$ws.Range("H134").Value = 4259
$ws.Range("I134").Value = 4259
$ws.Range("K134").Value = 12777
$ws.Range("M134").Value = -10242

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 140.2
$ws.Range("I7").Value = 100.25
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 100.25
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = 12.75
$ws.Range("N7").Value = -526

$ws.Range("H10").Value = 4003
$ws.Range("I10").Value = 4003
$ws.Range("K10").Value = 4003
$ws.Range("M10").Value = -3864

$ws.Range("H31").Value = 5266.577
$ws.Range("I31").Value = 3271.1667
$ws.Range("J31").Value = 6976.9287
$ws.Range("K31").Value = 3271.1667
$ws.Range("L31").Value = 6976.9287
$ws.Range("M31").Value = -2976.1667
$ws.Range("N31").Value = -7566.9287

$ws.Range("H34").Value = 5266.577
$ws.Range("I34").Value = 3271.1667
$ws.Range("J34").Value = 6976.9287
$ws.Range("K34").Value = 3271.1667
$ws.Range("L34").Value = 6976.9287
$ws.Range("M34").Value = -3069.1667
$ws.Range("N34").Value = -7380.9287

$ws.Range("H58").Value = 2389.5
$ws.Range("J58").Value = 2387.5
$ws.Range("L58").Value = 2387.5
$ws.Range("N58").Value = -2793.5

$ws.Range("H133").Value = 124600
$ws.Range("J133").Value = 124600
$ws.Range("L133").Value = 124600
$ws.Range("N133").Value = -129660

$ws.Range("H134").Value = 9413.429
$ws.Range("I134").Value = 10665.667
$ws.Range("J134").Value = 1900
$ws.Range("K134").Value = 31997.001
$ws.Range("L134").Value = 5700
$ws.Range("M134").Value = -29462.001
$ws.Range("N134").Value = -10770

$ws.Range("H136").Value = 2389.5
$ws.Range("J136").Value = 2387.5
$ws.Range("L136").Value = 7162.5
$ws.Range("N136").Value = -12262.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 199.85715
$ws.Range("I7").Value = 99.75
$ws.Range("J7").Value = 333.33334
$ws.Range("K7").Value = 299.25
$ws.Range("L7").Value = 1000.00002
$ws.Range("M7").Value = -187.25
$ws.Range("N7").Value = -1224.00002

$ws.Range("H68").Value = 1349.5
$ws.Range("J68").Value = 1499
$ws.Range("L68").Value = 4497
$ws.Range("N68").Value = -6119

$ws.Range("H71").Value = 1349.5
$ws.Range("J71").Value = 1499
$ws.Range("L71").Value = 13491
$ws.Range("N71").Value = -21603

$ws.Range("H107").Value = 2598.3635
$ws.Range("I107").Value = 3920.1667
$ws.Range("J107").Value = 1012.2
$ws.Range("K107").Value = 11760.5001
$ws.Range("L107").Value = 3036.6
$ws.Range("M107").Value = -9840.500100000001
$ws.Range("N107").Value = -6876.6

$ws.Range("H113").Value = 624.4
$ws.Range("I113").Value = 515.8333
$ws.Range("J113").Value = 787.25
$ws.Range("K113").Value = 1547.4999
$ws.Range("L113").Value = 2361.75
$ws.Range("M113").Value = 622.5001
$ws.Range("N113").Value = -6701.75

$ws.Range("H140").Value = 1399.2
$ws.Range("I140").Value = 1399.2
$ws.Range("K140").Value = 4197.6
$ws.Range("M140").Value = 982.3999999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 66666
$ws.Range("J45").Value = 66666
$ws.Range("L45").Value = 66666
$ws.Range("N45").Value = -67784

$ws.Range("H46").Value = 39496
$ws.Range("J46").Value = 39496
$ws.Range("L46").Value = 39496
$ws.Range("N46").Value = -39808

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents() | Out-Null

$ws.Range("H126").Value = 2034.7693
$ws.Range("I126").Value = 1690.4
$ws.Range("K126").Value = 5071.200000000001
$ws.Range("M126").Value = -2601.200000000001

$ws.Range("H132").Value = 4934.952
$ws.Range("I132").Value = 4850.5386
$ws.Range("J132").Value = 5072.125
$ws.Range("K132").Value = 14551.6158
$ws.Range("L132").Value = 15216.375
$ws.Range("M132").Value = -12021.6158
$ws.Range("N132").Value = -20276.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3416.1667
$ws.Range("J7").Value = 4500
$ws.Range("L7").Value = 4500
$ws.Range("N7").Value = -4724

$ws.Range("H122").Value = 3945
$ws.Range("I122").Value = 3945
$ws.Range("K122").Value = 11835
$ws.Range("M122").Value = -9385

$ws.Range("H126").Value = 3416.1667
$ws.Range("J126").Value = 4500
$ws.Range("L126").Value = 13500
$ws.Range("N126").Value = -18440

$ws.Range("H132").Value = 5473
$ws.Range("I132").Value = 4168.5
$ws.Range("J132").Value = 5994.8
$ws.Range("K132").Value = 12505.5
$ws.Range("L132").Value = 17984.4
$ws.Range("M132").Value = -9975.5
$ws.Range("N132").Value = -23044.4

$ws.Range("H136").Value = 24451.092
$ws.Range("I136").Value = 2881.5557
$ws.Range("K136").Value = 8644.667099999999
$ws.Range("M136").Value = -6094.667099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 38906.668
$ws.Range("J54").Value = 38906.668
$ws.Range("L54").Value = 38906.668
$ws.Range("N54").Value = -39946.668

$ws.Range("H132").Value = 2695.7
$ws.Range("I132").Value = 2399.2693
$ws.Range("J132").Value = 4622.5
$ws.Range("K132").Value = 7197.8079
$ws.Range("L132").Value = 13867.5
$ws.Range("M132").Value = -4667.8079
$ws.Range("N132").Value = -18927.5
